$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '27.289.25'
$ws.Cells.Item(2,5).Value = '  -0.90%  '

# Row 3
$ws.Cells.Item(3,4).Value = '1.783.83'
$ws.Cells.Item(3,5).Value = '  -2.17%  '

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '1.004'
$ws.Cells.Item(4,5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '340.77'
$ws.Cells.Item(5,5).Value = '  -0.65%  '

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.3960'
$ws.Cells.Item(7,5).Value = '  +3.65%  '

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.3456'
$ws.Cells.Item(8,5).Value = '  -2.25%  '

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '47.97'
$ws.Cells.Item(9,5).Value = '  -3.97%  '

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '1.195'
$ws.Cells.Item(10,5).Value = '  -3.43%  '

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07467'
$ws.Cells.Item(11,5).Value = '  -3.55%  '

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '1.001'
$ws.Cells.Item(12,5).Value = '  -0.15%  '

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '21.73'
$ws.Cells.Item(13,5).Value = '  -2.29%  '

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '6.482'
$ws.Cells.Item(14,5).Value = '  -2.05%  '

# Row 15
$ws.Cells.Item(15,4).Value = '1.784.37'
$ws.Cells.Item(15,5).Value = '  -2.23%  '

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '7.109'
$ws.Cells.Item(16,5).Value = '  -1.56%  '

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.00001093'
$ws.Cells.Item(17,5).Value = '  -2.85%  '

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '0.06699'
$ws.Cells.Item(18,5).Value = '  -0.59%  '

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '84.29'
$ws.Cells.Item(19,5).Value = '  -3.13%  '

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '1.001'
$ws.Cells.Item(20,5).Value = '  -0.07%  '

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '17.64'
$ws.Cells.Item(21,5).Value = '  +0.29%  '

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '6.513'
$ws.Cells.Item(22,5).Value = '  -0.39%  '

# Row 23
$ws.Cells.Item(23,4).Value = '27.290.06'
$ws.Cells.Item(23,5).Value = '  -0.87%  '

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '12.42'
$ws.Cells.Item(24,5).Value = '  -5.72%  '

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '2.380'
$ws.Cells.Item(25,5).Value = '  -4.12%  '

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '1.475'
$ws.Cells.Item(26,5).Value = '  -0.39%  '

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '21.17'
$ws.Cells.Item(27,5).Value = '  -3.76%  '

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '2.497'
$ws.Cells.Item(28,5).Value = '  -6.94%  '

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '157.52'
$ws.Cells.Item(29,5).Value = '  +3.08%  '

# Row 30
$ws.Cells.Item(30,4).Value = '1.986.50'
$ws.Cells.Item(30,5).Value = '  -2.14%  '

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '136.48'
$ws.Cells.Item(31,5).Value = '  +0.65%  '

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '4.026'
$ws.Cells.Item(32,5).Value = '  -1.65%  '

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '5.963'
$ws.Cells.Item(33,5).Value = '  -5.78%  '

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '0.08835'
$ws.Cells.Item(34,5).Value = '  +0.57%  '

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '13.01'
$ws.Cells.Item(35,5).Value = '  -6.41%  '

# Row 36
$ws.Cells.Item(36,2).Value = 'VeChain'
$ws.Cells.Item(36,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.02433'
$ws.Cells.Item(36,5).Value = '  +1.28%  '

# Row 37
$ws.Cells.Item(37,2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(37,3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '1.620'
$ws.Cells.Item(37,5).Value = '  -4.55%  '

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '5.414'
$ws.Cells.Item(38,5).Value = '  -3.77%  '

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '0.06461'
$ws.Cells.Item(39,5).Value = '  -0.90%  '

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.6819'
$ws.Cells.Item(40,5).Value = '  -2.82%  '

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.2204'
$ws.Cells.Item(41,5).Value = '  -2.49%  '

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '1.252'
$ws.Cells.Item(42,5).Value = '  -4.00%  '

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '8.389'
$ws.Cells.Item(43,5).Value = '  -7.99%  '

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '14.40'
$ws.Cells.Item(44,5).Value = '  -2.75%  '

# Row 45
$ws.Cells.Item(45,5).Value = '  -0.08%  '

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.6394'
$ws.Cells.Item(46,5).Value = '  -3.40%  '

# Row 47
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '3.879'
$ws.Cells.Item(47,5).Value = '  -1.66%  '

# Row 48
$ws.Cells.Item(48,2).Value = 'NEARProtocol'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '2.137'
$ws.Cells.Item(48,5).Value = '  -2.40%  '

# Row 49
$ws.Cells.Item(49,2).Value = 'Quant'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '132.45'
$ws.Cells.Item(49,5).Value = '  -0.59%  '

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.07137'
$ws.Cells.Item(50,5).Value = '  -2.40%  '

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '79.29'
$ws.Cells.Item(51,5).Value = '  -2.29%  '
